# Latest Updates from Costco
$wb = $excel.ActiveWorkbook

# Rename the "Aboutus" sheet to "AboutUs"
$ws1 = $wb.Worksheets.Item("Aboutus")
$ws1.Name = "AboutUs"

$ws2 = $wb.Worksheets.Item("Membership")

# Insert a new row for "Kirkland Signature" right after "Jobs" (row 6),
# pushing the rest of the list down by one.
$ws1.Rows.Item(7).Insert()
$ws1.Range("A7").Value = "Kirkland Signature"

# Fix the ampersand in the recipe videos link text.
$ws1.Range("A11").Value = "Quick & Easy Recipe Videos"

# Update the selections to match the latest edits on each sheet.
$ws1.Range("A11").Select()
$ws2.Range("A7").Select()

$ws1.Activate()
